# Add a "monthly amount paid" tracker to Sheet1: new header columns G:S,
# a running Total-paid formula down to row 10, a couple of extra data rows,
# and a couple of new columns' widths. Also normalize Sheet2/Sheet3 a touch.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New header row (G1:S1) ---------------------------------------------
$headerCols = @("G","H","I","J","K","L","M","N","O","P","Q","R","S")
$headerVals = @("Amount Paid","Baisakh","Jestha","Ashad","Shrawan","Bhadra","Karktik","Mangsir","Poush","Magh","Falgun","Chaitra","Total paid")
for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $headerVals[$i]
}

# --- Row 2: Amount-paid figure for existing first record ----------------
$ws.Range("H2").Value = 1

# --- Row 6: new record (Sandeep Thapa) -----------------------------------
$ws.Range("A6").Value = "Sandeep Thapa"
$ws.Range("B6").Value = "Babu Ram Thapa"
$ws.Range("C6").Value = "Shree Thapa"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "21"
$ws.Range("E6").Value = "chitwan"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "989911"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "2000"
$ws.Range("H6").Value = 2000
$ws.Range("I6").Value = 2

# --- S column running totals, rows 2-10 ----------------------------------
for ($r = 2; $r -le 10; $r++) {
    $formula = "=SUM(Q$r,H$r,I$r,J$r,K$r,L$r" + ":M$r,N$r,O$r,P$r,R$r)"
    $ws.Range("S$r").Formula = $formula
}

# --- Row 11: new record (Sandeep Thapa again, pasted as plain text) ------
$ws.Range("A11").Value = "Sandeep Thapa"
$ws.Range("B11").Value = "Babu Ram Thapa"
$ws.Range("C11").Value = "Shree Thapa"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21"
$ws.Range("E11").Value = "chitwan"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "989911"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "2000"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "2000"

# --- Row 12: new record (test row) ----------------------------------------
$ws.Range("A12").Value = "drter"
$ws.Range("B12").Value = "rtyt"
$ws.Range("C12").Value = "ggh"
$ws.Range("D12").Value = "hghk"
$ws.Range("E12").Value = "huhu"
$ws.Range("F12").Value = "guhj"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "200"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "200"

# --- Column widths for the new columns ------------------------------------
$ws.Columns("F:G").ColumnWidth = 14.88671875
$ws.Columns("H:H").ColumnWidth = 11.6640625

# --- View: scroll right a bit and select G6, matching the saved view -----
[void]$ws.Range("G6").Select()

# --- Page setup: portrait orientation -------------------------------------
$ws.PageSetup.Orientation = 1

# --- Sheet2 / Sheet3: touch the default row height (cosmetic resave) -----
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")
